$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.535.74"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "3.419.06"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.93%  "
$ws.Range("E7").Value = "  +6.41%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "3.424.51"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "4.012.10"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("E15").Value = "  -3.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.12%  "
$ws.Range("D17").Value = "64.577.43"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "3.424.04"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.556"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("E26").Value = "  -4.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.13%  "
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  +4.81%  "
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("E35").Value = "  +6.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("D41").Value = "2.848.23"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.30%  "
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.770"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "320.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.93%  "
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("E51").Value = "  -1.98%  "
